# Update Bank Deposit data - 2025-12-23T10:14:24.863Z
#
# - Row 13: the Cash deposit collected on 02-12-2025 was actually banked on
#   03-12-2025 (Deposit Date) for Rs. 8770 (Deposit Amount).
# - Five new rows (15-19) record the 03-12-2025 collection being deposited.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that hold date-like text ("A", "D", "G") must be pre-formatted as
# Text so Excel keeps the literal string instead of re-interpreting it as a
# serial date number.
$dateTextFormat = "@"

# --- Row 13: fill in the deposit date / amount now that it is known -------
$ws.Range("D13").NumberFormat = $dateTextFormat
$ws.Range("D13").Value = "2025-12-03"
$ws.Range("E13").Value = 8770

# --- New rows 15-19: the 03-12-2025 collection deposits -------------------
$newRows = @(
    @{ Row = 15; A = "03-12-2025"; B = "010965012-Medha Sub Division Office Coll."; E = 23120 },
    @{ Row = 16; A = "03-12-2025"; B = "020965017-Kai Lalsingrao Shinde Gr.Bid.S.S.Pat.Ltd Kudal Br. Kudal"; E = 78770 },
    @{ Row = 17; A = "03-12-2025"; B = "020965018-Kai Lalsingrao Shinde Gr.Big.Sheti Sah.Pat.Ltd. Br. Medha"; E = 64270 },
    @{ Row = 18; A = "03-12-2025"; B = "020965019-SHRI DATTATRAY MAHARAJ KALAMBE SAH. PAT. LTD.DAPAWADI"; E = 16600 },
    @{ Row = 19; A = "03-12-2025"; B = "020965020-KAI.LALSINGRAO BAPUSO SHINDE SAH.PAT.LTD.,KUDAL BR.SAYGAON"; E = 11350 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    $ws.Range("A$rowNum").NumberFormat = $dateTextFormat
    $ws.Range("A$rowNum").Value = $r.A

    $ws.Range("B$rowNum").Value = $r.B

    $ws.Range("C$rowNum").Value = "Cash"

    $ws.Range("D$rowNum").NumberFormat = $dateTextFormat
    $ws.Range("D$rowNum").Value = "2025-12-03"

    $ws.Range("E$rowNum").Value = $r.E

    $ws.Range("F$rowNum").Value = ""

    $ws.Range("G$rowNum").NumberFormat = $dateTextFormat
    $ws.Range("G$rowNum").Value = "2025-12-23"
}

Write-Host "Bank Deposit data updated"
